$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to their "native" (PascalCase) equivalents.
$ws.Range("A1").Value = "StandardCategoryID"
$ws.Range("B1").Value = "DetailedStandardCategoryID"

# Match the "General" number format used by the other headers (D1/E1)
# by copying D1's formatting onto the renamed cells, rather than assigning
# NumberFormat directly (which would create a brand new style entry).
$ws.Range("D1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# Reset the clipboard/marching ants after the paste.
$excel.CutCopyMode = $false

# Put the selection back on A1 (the default/top-left cell).
$ws.Range("A1").Select()
